# Append/update: set the "取得日時" (acquired datetime) column A values
# for data rows 2-13 on the "ランサーズ" sheet to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-29 01:48:15"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
